# Weekly refresh: insert the newest day's record at the top (row 2) and
# push every existing record down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the first data row (row 2). Excel's default
# Insert() copies formatting from the row above (the bold header row), so
# strip that back to Normal and only re-apply the date number format that
# every other row in column D carries.
$ws.Rows.Item(2).Insert()
$ws.Range("A2:R2").Style = "Normal"
$ws.Range("D2").NumberFormat = $ws.Range("D3").NumberFormat

# Populate the new top row with the latest market reading.
$ws.Range("A2").Value = 11
$ws.Range("B2").Value = "Vega Monumental Concepción"
$ws.Range("C2").Value = "Bíobío"
$ws.Range("D2").Value = 44860
$ws.Range("E2").Value = 8
$ws.Range("F2").Value = 100112037
$ws.Range("G2").Value = "Cebollín"
$ws.Range("H2").Value = "Sin especificar"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 250
$ws.Range("K2").Value = 2000
$ws.Range("L2").Value = 2200
$ws.Range("M2").Value = 2080
$ws.Range("N2").Value = "$/paquete 36 unidades"
$ws.Range("O2").Value = "Región Metropolitana"
$ws.Range("P2").Value = 58
$ws.Range("Q2").Value = 36
$ws.Range("R2").Value = "Hortaliza"
